# Refresh cryptos list (prices / 1h volumes) per GitHub Actions update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces text entry, matching the original sheet where
# every value in columns B:E is stored as text (prevents Excel from
# auto-coercing numeric-looking strings like '0.506' or '19.67' into reals).
# Re-applying the Normal style afterwards clears the quote-prefix marker
# Excel would otherwise leave on the cell, keeping formatting untouched.
$ws.Range('D2').Value = "'25.810.41"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -0.30%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('E3').Value = "'  -0.18%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  -0.23%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('E5').Value = "'  +0.08%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'0.506"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -0.72%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'  -0.21%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('E9').Value = "'  -0.87%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'19.67"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E11').Value = "'  +1.59%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = "'  -0.11%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('B13').Value = "'WrappedliquidstakedEther2.0"
$ws.Range('B13').Style = 'Normal'
$ws.Range('C13').Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range('C13').Style = 'Normal'
$ws.Range('D13').Value = "'1.860.36"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -0.17%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('B14').Value = "'WrappedEther"
$ws.Range('B14').Style = 'Normal'
$ws.Range('C14').Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range('C14').Style = 'Normal'
$ws.Range('D14').Value = "'1.634.93"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -1.36%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'0.560"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +0.02%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('E16').Value = "'  +0.08%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'62.86"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  -0.67%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'25.821.07"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E19').Value = "'  -0.23%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = "'  +1.68%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'194.31"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  +0.36%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'9.95"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  -0.16%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'6.28"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +2.05%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('E24').Value = "'  -0.18%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = "'  +2.27%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'142.75"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  +3.38%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = "'  +0.00%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'6.91"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  +0.94%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'  +0.00%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = "'  -0.34%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = "'  -0.44%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = "'  +1.73%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = "'  -0.24%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = "'  +0.54%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'2.39"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  -0.11%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'0.904"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -0.10%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'1.135.24"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  -0.24%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = "'  -1.60%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'  -1.65%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Value = "'  -0.68%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').Value = "'  +0.39%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'5.59"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +2.08%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'100.58"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +0.82%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('E44').Value = "'  +0.56%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'1.769.98"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  -0.36%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('E46').Value = "'  +1.14%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'55.26"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  -0.57%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = "'  -0.43%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'0.417"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -2.34%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('B50').Value = "'EnergySwap"
$ws.Range('B50').Style = 'Normal'
$ws.Range('C50').Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range('C50').Style = 'Normal'
$ws.Range('D50').Value = "'7.55"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  -2.99%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('B51').Value = "'Frax"
$ws.Range('B51').Style = 'Normal'
$ws.Range('C51').Value = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range('C51').Style = 'Normal'
$ws.Range('D51').Value = "'1.00"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +0.22%  "
$ws.Range('E51').Style = 'Normal'
